$wb = $excel.ActiveWorkbook

# --- Overview sheet: update status text for the 02a80f9e row ---
$wsOverview = $wb.Worksheets.Item("Overview")
$wsOverview.Range("B3").Value = "Ready for handoff"
$wsOverview.Range("C3").Value = "Ready for handoff"

# --- zh-cn sheet: update status + latest handoff datetime for the 02a80f9e row ---
$wsZhCn = $wb.Worksheets.Item("zh-cn")
$wsZhCn.Range("B3").Value = "Ready for handoff"
$wsZhCn.Range("D3").Value = "2016-01-15 07:54:39"

# --- de-de sheet: update status + latest handoff datetime for the 02a80f9e row ---
$wsDeDe = $wb.Worksheets.Item("de-de")
$wsDeDe.Range("B3").Value = "Ready for handoff"
$wsDeDe.Range("D3").Value = "2016-01-15 07:54:50"
